$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look like plain numbers (e.g. "246.17", "1.0000",
# "0.04620"). Force those specific cells to text format first so Excel
# keeps the exact original string (incl. trailing zeros) instead of
# silently converting them to numeric values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.616.15'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '1.739.19'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '246.17'
$ws.Range("E5").Value = '  +0.14%  '

$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").Value = '0.4926'
$ws.Range("E7").Value = '  +2.68%  '

$ws.Range("E8").Value = '  -0.55%  '

$ws.Range("D9").Value = '0.06276'
$ws.Range("E9").Value = '  +0.82%  '

$ws.Range("D10").Value = '1.749.02'
$ws.Range("E10").Value = '  +1.21%  '

$ws.Range("D11").Value = '0.07044'
$ws.Range("E11").Value = '  -1.26%  '

$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").Value = '0.6138'
$ws.Range("E13").Value = '  -0.79%  '

$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").Value = '77.97'
$ws.Range("E15").Value = '  +0.95%  '

$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '26.635.73'
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").Value = '1.0000'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").Value = '0.000007265'
$ws.Range("E19").Value = '  +4.38%  '

$ws.Range("E20").Value = '  -1.25%  '

$ws.Range("D21").Value = '1.973.25'
$ws.Range("E21").Value = '  +1.09%  '

$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").Value = '8.711'
$ws.Range("E23").Value = '  -2.68%  '

$ws.Range("D24").Value = '5.273'
$ws.Range("E24").Value = '  -0.52%  '

$ws.Range("D25").Value = '139.04'
$ws.Range("E25").Value = '  +1.87%  '

$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("E27").Value = '  +1.01%  '

$ws.Range("D29").Value = '107.37'
$ws.Range("E29").Value = '  +0.50%  '

$ws.Range("D30").Value = '4.028'
$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("D31").Value = '0.08047'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '3.728'
$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").Value = '0.04620'
$ws.Range("E33").Value = '  +1.11%  '

$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '0.9994'
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.613'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.015'
$ws.Range("E36").Value = '  +2.53%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6380'
$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.066'
$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.9039'
$ws.Range("E39").Value = '  -3.31%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.426'
$ws.Range("E40").Value = '  +0.53%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01505'
$ws.Range("E42").Value = '  +0.34%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '101.98'
$ws.Range("E43").Value = '  -3.77%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.432'
$ws.Range("E44").Value = '  -4.89%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3932'
$ws.Range("E45").Value = '  +0.49%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.864'
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1186'
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05389'
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.61'
$ws.Range("E49").Value = '  -1.53%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.810'
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.255'
$ws.Range("E51").Value = '  -0.95%  '
